# Update "想去人数" (F) and "最低票价" (G) values for matching events across
# the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row -> new F value (and optionally new G value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13573
$ws1.Range("F6").Value  = 26
$ws1.Range("F10").Value = 83
$ws1.Range("F13").Value = 13578
$ws1.Range("G13").Value = 49.9
$ws1.Range("F15").Value = 603
$ws1.Range("F16").Value = 8969
$ws1.Range("F18").Value = 8064
$ws1.Range("F19").Value = 256
$ws1.Range("F26").Value = 1022
$ws1.Range("F30").Value = 209
$ws1.Range("F31").Value = 191

# Sheet "全部类型" - same events, rows shifted by +2 for the last two (30,31 -> 32,33)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 13573
$ws4.Range("F6").Value  = 26
$ws4.Range("F10").Value = 83
$ws4.Range("F13").Value = 13578
$ws4.Range("G13").Value = 49.9
$ws4.Range("F15").Value = 603
$ws4.Range("F16").Value = 8969
$ws4.Range("F18").Value = 8064
$ws4.Range("F19").Value = 256
$ws4.Range("F26").Value = 1022
$ws4.Range("F32").Value = 209
$ws4.Range("F33").Value = 191
